# Update the asset tags list: replace the old (36-item) CAPEX/sector tag
# list with the new (17-item) infrastructure asset-class list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Asset List Level 1",
    "Aviation",
    "Bridges",
    "Broadband",
    "Dams",
    "Drinking Water",
    "Energy",
    "Hazardous Waste",
    "Inland Waterways",
    "Ports",
    "Rail",
    "Roads",
    "Schools",
    "Solid Waste",
    "Stormwater",
    "Transit",
    "Waste Water"
)

# Drop the old rows that no longer exist in the new list (18-36), shrinking
# the sheet down to A1:A17.
$ws.Range("A18:A36").EntireRow.Delete()

# Overwrite A1:A17 with the new asset tag values (A1 keeps its existing bold
# header style).
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Move the active selection to the first empty row below the list, matching
# the saved view state.
$ws.Range("A18").Select() | Out-Null
